$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "东百集团"
$ws.Range("B2").Value = "平潭发展"
$ws.Range("C2").Value = "平潭发展"
$ws.Range("A3").Value = "航天发展"
$ws.Range("B3").Value = "永辉超市"
$ws.Range("C3").Value = "东百集团"
$ws.Range("A4").Value = "永辉超市"
$ws.Range("B4").Value = "东百集团"
$ws.Range("C4").Value = "合富中国"
$ws.Range("A5").Value = "平潭发展"
$ws.Range("B5").Value = "雪人集团"
$ws.Range("C5").Value = "航天发展"
$ws.Range("A6").Value = "雪人集团"
$ws.Range("B6").Value = "航天发展"
$ws.Range("C6").Value = "永辉超市"
$ws.Range("A7").Value = "浙江世宝"
$ws.Range("B7").Value = "西部材料"
$ws.Range("C7").Value = "雪人集团"
$ws.Range("A8").Value = "西部材料"
$ws.Range("B8").Value = "百大集团"
$ws.Range("C8").Value = "再升科技"
$ws.Range("A9").Value = "顺灏股份"
$ws.Range("B9").Value = "美年健康"
$ws.Range("C9").Value = "西部材料"
$ws.Range("A10").Value = "航天电子"
$ws.Range("B10").Value = "顺灏股份"
$ws.Range("C10").Value = "浙江世宝"
$ws.Range("A11").Value = "百大集团"
$ws.Range("B11").Value = "山子高科"
$ws.Range("C11").Value = "山子高科"
$ws.Range("A12").Value = "合富中国"
$ws.Range("B12").Value = "合富中国"
$ws.Range("C12").Value = "百大集团"
$ws.Range("A13").Value = "航天动力"
$ws.Range("B13").Value = "浙江世宝"
$ws.Range("C13").Value = "海峡创新"
$ws.Range("A14").Value = "山子高科"
$ws.Range("B14").Value = "中国中免"
$ws.Range("C14").Value = "航天电子"
$ws.Range("A15").Value = "美年健康"
$ws.Range("B15").Value = "航天电子"
$ws.Range("C15").Value = "金圆股份"
$ws.Range("A16").Value = "王子新材"
$ws.Range("B16").Value = "中央商场"
$ws.Range("C16").Value = "九牧王"
$ws.Range("A17").Value = "再升科技"
$ws.Range("B17").Value = "中国卫星"
$ws.Range("C17").Value = "顺灏股份"
$ws.Range("A18").Value = "中国中免"
$ws.Range("B18").Value = "海南发展"
$ws.Range("C18").Value = "海南发展"
$ws.Range("A19").Value = "海南发展"
$ws.Range("B19").Value = "再升科技"
$ws.Range("C19").Value = "中国中免"
$ws.Range("A20").Value = "中央商场"
$ws.Range("B20").Value = "航天动力"
$ws.Range("C20").Value = "神剑股份"
$ws.Range("A21").Value = "鹭燕医药"
$ws.Range("B21").Value = "C沐曦-U"
$ws.Range("C21").Value = "中央商场"
